$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-7 from 45208 (2023-10-09)
# to 45212 (2023-10-13), preserving existing number formatting/style.
foreach ($row in 2..7) {
    $ws.Range("C$row").Value = 45212
}
